$d = $word.ActiveDocument

# --- Paragraph 1: "Dispute Documentation Guide Template" (Bold, sz32)
#     -> "Dispute Documentation Guide" styled as Heading1, with no direct
#     character formatting left on the run.
$p1 = $d.Paragraphs(1)
$p1.Range.Delete()
$d.Range(0, 0).InsertBefore("Dispute Documentation Guide`r")
$d.Paragraphs(1).Style = "Heading1"

# --- Paragraph 2: "Date: 2026-01-21" -> "Tenant: {{TENANT_NAME}}"
$d.Content.Find.Execute("Date: 2026-01-21", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Tenant: {{TENANT_NAME}}", 2)

# --- Paragraph 3: "Not legal advice. General-purpose template." + <w:br/>
#     -> "Property: {{PROPERTY_ADDRESS}}" (drop the line break run content)
$p3 = $d.Paragraphs(3)
$p3.Range.Delete()
$d.Paragraphs(3).Range.InsertBefore("Property: {{PROPERTY_ADDRESS}}`r")

# --- Paragraph 4: "Timeline" -> "Issue Summary: {{ISSUE_SUMMARY}}"
$d.Content.Find.Execute("Timeline", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Issue Summary: {{ISSUE_SUMMARY}}", 2)

# --- Paragraph 5: "Evidence" -> "Timeline:"
$d.Content.Find.Execute("Evidence", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Timeline:", 2)

# --- Paragraph 6: "Resolution" -> "- {{DATE}}: {{EVENT}}"
$d.Content.Find.Execute("Resolution", $true, $false, $false, $false, `
    $false, $true, 1, $false, "- {{DATE}}: {{EVENT}}", 2)

# --- Paragraph 7: "File Naming" -> "- {{DATE}}: {{EVENT}}"
$d.Content.Find.Execute("File Naming", $true, $false, $false, $false, `
    $false, $true, 1, $false, "- {{DATE}}: {{EVENT}}", 2)

# --- New paragraphs 8 & 9 appended at the end of the document
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$pLast.Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = "Supporting Evidence:"

$pLast2 = $d.Paragraphs($d.Paragraphs.Count)
$pLast2.Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = "- {{EVIDENCE_ITEM}}"
